# Final settlement match logic added
# Renumber the M0xx match groups (rows 11-28) and mark the later groups
# (now M025-M029) as 'Settlement' matches with updated audit info text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group that becomes M030 (rows 11-13), was M025. Match type stays 'Manual'. ---
$ws.Range("A11").Value = "M030"
$ws.Range("A12").Value = "M030"
$ws.Range("A13").Value = "M030"

# --- Group that becomes M025 (rows 14-16), was M026. Now a 'Settlement' match. ---
$ws.Range("A14").Value = "M025"
$ws.Range("L14").Value = "Settlement"

$ws.Range("A15").Value = "M025"
$ws.Range("B15").Value = "Settlement Match (ID: 11370) - 'Final Settlement' keyword found`nLender Amount: 99317.00`nBorrower Amount: 99317.00"
$ws.Range("L15").Value = "Settlement"

$ws.Range("A16").Value = "M025"
$ws.Range("L16").Value = "Settlement"

# --- Group that becomes M026 (rows 17-19), was M027. Now a 'Settlement' match. ---
$ws.Range("A17").Value = "M026"
$ws.Range("L17").Value = "Settlement"

$ws.Range("A18").Value = "M026"
$ws.Range("B18").Value = "Settlement Match (ID: 10199) - 'Final Settlement' keyword found`nLender Amount: 20130.00`nBorrower Amount: 20130.00"
$ws.Range("L18").Value = "Settlement"

$ws.Range("A19").Value = "M026"
$ws.Range("L19").Value = "Settlement"

# --- Group that becomes M027 (rows 20-22), was M028. Now a 'Settlement' match. ---
$ws.Range("A20").Value = "M027"
$ws.Range("L20").Value = "Settlement"

$ws.Range("A21").Value = "M027"
$ws.Range("B21").Value = "Settlement Match (ID: 11711) - 'Final Settlement' keyword found`nLender Amount: 94109.00`nBorrower Amount: 94109.00"
$ws.Range("L21").Value = "Settlement"

$ws.Range("A22").Value = "M027"
$ws.Range("L22").Value = "Settlement"

# --- Group that becomes M028 (rows 23-25), was M029. Now a 'Settlement' match. ---
$ws.Range("A23").Value = "M028"
$ws.Range("L23").Value = "Settlement"

$ws.Range("A24").Value = "M028"
$ws.Range("B24").Value = "Settlement Match (ID: 11134) - 'Final Settlement' keyword found`nLender Amount: 13909.00`nBorrower Amount: 13909.00"
$ws.Range("L24").Value = "Settlement"

$ws.Range("A25").Value = "M028"
$ws.Range("L25").Value = "Settlement"

# --- Group that becomes M029 (rows 26-28), was M030. Now a 'Settlement' match. ---
$ws.Range("A26").Value = "M029"
$ws.Range("L26").Value = "Settlement"

$ws.Range("A27").Value = "M029"
$ws.Range("B27").Value = "Settlement Match (ID: 12107) - 'Final Settlement' keyword found`nLender Amount: 93314.00`nBorrower Amount: 93314.00"
$ws.Range("L27").Value = "Settlement"

$ws.Range("A28").Value = "M029"
$ws.Range("L28").Value = "Settlement"
